$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bordered (blank) cells in columns C:F for rows 1-3 ---
$ws.Range("C1:F3").Borders.LineStyle = 1

# --- Row 3: give A3 and B3 a full border (they had none before) ---
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("B3").Borders.LineStyle = 1

# --- New bordered (blank) cells in column F for rows 6-15 ---
$ws.Range("F6:F15").Borders.LineStyle = 1

# --- Merge E6:E16 (like the other data columns A/B/C/D) ---
# Clear borders first so the merge operation (which always strips the
# inner borders of a merged block) doesn't have to fabricate new
# mixed-border combinations; then apply the final box border once merged.
$ws.Range("E6:E16").Borders.LineStyle = -4142
$ws.Range("E6:E16").Merge()
$ws.Range("E6:E16").Borders.LineStyle = 1

# --- E6 becomes the "signature" label, styled like A6:D6 (bordered, centered, wrapped) ---
$ws.Range("E6").Value = "signature"
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").VerticalAlignment = -4108
$ws.Range("E6").WrapText = $true

# --- E16 should match the rest of row 16 (bottom border only) ---
$ws.Range("E16").Borders.LineStyle = -4142
$ws.Range("E16").Borders.Item(9).LineStyle = 1
